# Update the natmi LR-pair recalculated TPM output (Clcf1-Crlf1).
# The source data now covers the full 4x4 sending/target cluster grid
# (ECs, FAPs, MuSCs, Resolving-Mac), growing the table from 12 to 16 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Clcf1"
$ws.Range("C2").Value = "Crlf1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.874436
$ws.Range("H2").Value = 5.623308
$ws.Range("I2").Value = 0.1442186763702422
$ws.Range("J2").Value = 0.1442186763702422
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3711636666666667
$ws.Range("N2").Value = 1.113491
$ws.Range("O2").Value = 0.0322568113697387
$ws.Range("P2").Value = 0.03225681136973869
$ws.Range("Q2").Value = 0.695722538692
$ws.Range("R2").Value = 6.261502848228
$ws.Range("S2").Value = 0.004652034639668294
$ws.Range("T2").Value = 0.004652034639668292

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Clcf1"
$ws.Range("C3").Value = "Crlf1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.874436
$ws.Range("H3").Value = 5.623308
$ws.Range("I3").Value = 0.1442186763702422
$ws.Range("J3").Value = 0.1442186763702422
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.123096
$ws.Range("N3").Value = 18.369288
$ws.Range("O3").Value = 0.5321413985496108
$ws.Range("P3").Value = 0.5321413985496107
$ws.Range("Q3").Value = 11.477351573856
$ws.Range("R3").Value = 103.296164164704
$ws.Range("S3").Value = 0.07674472814063439
$ws.Range("T3").Value = 0.07674472814063436

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Clcf1"
$ws.Range("C4").Value = "Crlf1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.874436
$ws.Range("H4").Value = 5.623308
$ws.Range("I4").Value = 0.1442186763702422
$ws.Range("J4").Value = 0.1442186763702422
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.943314666666667
$ws.Range("N4").Value = 14.829944
$ws.Range("O4").Value = 0.4296098542617661
$ws.Range("P4").Value = 0.4296098542617661
$ws.Range("Q4").Value = 9.265926970528
$ws.Range("R4").Value = 83.393342734752
$ws.Range("S4").Value = 0.06195776453724456
$ws.Range("T4").Value = 0.06195776453724454

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Clcf1"
$ws.Range("C5").Value = "Crlf1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.874436
$ws.Range("H5").Value = 5.623308
$ws.Range("I5").Value = 0.1442186763702422
$ws.Range("J5").Value = 0.1442186763702422
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06894633333333333
$ws.Range("N5").Value = 0.206839
$ws.Range("O5").Value = 0.005991935818884376
$ws.Range("P5").Value = 0.005991935818884375
$ws.Range("Q5").Value = 0.129235489268
$ws.Range("R5").Value = 1.163119403412
$ws.Range("S5").Value = 0.0008641490526949479
$ws.Range("T5").Value = 0.0008641490526949476

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Clcf1"
$ws.Range("C6").Value = "Crlf1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.197979
$ws.Range("H6").Value = 9.593937
$ws.Range("I6").Value = 0.2460517715407892
$ws.Range("J6").Value = 0.2460517715407892
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3711636666666667
$ws.Range("N6").Value = 1.113491
$ws.Range("O6").Value = 0.0322568113697387
$ws.Range("P6").Value = 0.03225681136973869
$ws.Range("Q6").Value = 1.186973611563
$ws.Range("R6").Value = 10.682762504067
$ws.Range("S6").Value = 0.007936845581781279
$ws.Range("T6").Value = 0.007936845581781276

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Clcf1"
$ws.Range("C7").Value = "Crlf1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.197979
$ws.Range("H7").Value = 9.593937
$ws.Range("I7").Value = 0.2460517715407892
$ws.Range("J7").Value = 0.2460517715407892
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.123096
$ws.Range("N7").Value = 18.369288
$ws.Range("O7").Value = 0.5321413985496108
$ws.Range("P7").Value = 0.5321413985496107
$ws.Range("Q7").Value = 19.581532422984
$ws.Range("R7").Value = 176.233791806856
$ws.Range("S7").Value = 0.1309343338233249
$ws.Range("T7").Value = 0.1309343338233249

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Clcf1"
$ws.Range("C8").Value = "Crlf1"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.197979
$ws.Range("H8").Value = 9.593937
$ws.Range("I8").Value = 0.2460517715407892
$ws.Range("J8").Value = 0.2460517715407892
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.943314666666667
$ws.Range("N8").Value = 14.829944
$ws.Range("O8").Value = 0.4296098542617661
$ws.Range("P8").Value = 0.4296098542617661
$ws.Range("Q8").Value = 15.808616494392
$ws.Range("R8").Value = 142.277548449528
$ws.Range("S8").Value = 0.1057062657124878
$ws.Range("T8").Value = 0.1057062657124878

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Clcf1"
$ws.Range("C9").Value = "Crlf1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.197979
$ws.Range("H9").Value = 9.593937
$ws.Range("I9").Value = 0.2460517715407892
$ws.Range("J9").Value = 0.2460517715407892
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.06894633333333333
$ws.Range("N9").Value = 0.206839
$ws.Range("O9").Value = 0.005991935818884376
$ws.Range("P9").Value = 0.005991935818884375
$ws.Range("Q9").Value = 0.220488926127
$ws.Range("R9").Value = 1.984400335143
$ws.Range("S9").Value = 0.00147432642319521
$ws.Range("T9").Value = 0.00147432642319521

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Clcf1"
$ws.Range("C10").Value = "Crlf1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.825289333333334
$ws.Range("H10").Value = 20.475868
$ws.Range("I10").Value = 0.5251361975000832
$ws.Range("J10").Value = 0.5251361975000832
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3711636666666667
$ws.Range("N10").Value = 1.113491
$ws.Range("O10").Value = 0.0322568113697387
$ws.Range("P10").Value = 0.03225681136973869
$ws.Range("Q10").Value = 2.533299415020889
$ws.Range("R10").Value = 22.799694735188
$ws.Range("S10").Value = 0.01693921926618203
$ws.Range("T10").Value = 0.01693921926618203

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Clcf1"
$ws.Range("C11").Value = "Crlf1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 6.825289333333334
$ws.Range("H11").Value = 20.475868
$ws.Range("I11").Value = 0.5251361975000832
$ws.Range("J11").Value = 0.5251361975000832
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 6.123096
$ws.Range("N11").Value = 18.369288
$ws.Range("O11").Value = 0.5321413985496108
$ws.Range("P11").Value = 0.5321413985496107
$ws.Range("Q11").Value = 41.79190181577601
$ws.Range("R11").Value = 376.127116341984
$ws.Range("S11").Value = 0.2794467105667189
$ws.Range("T11").Value = 0.2794467105667189

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Clcf1"
$ws.Range("C12").Value = "Crlf1"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 6.825289333333334
$ws.Range("H12").Value = 20.475868
$ws.Range("I12").Value = 0.5251361975000832
$ws.Range("J12").Value = 0.5251361975000832
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.943314666666667
$ws.Range("N12").Value = 14.829944
$ws.Range("O12").Value = 0.4296098542617661
$ws.Range("P12").Value = 0.4296098542617661
$ws.Range("Q12").Value = 33.73955286571022
$ws.Range("R12").Value = 303.6559757913921
$ws.Range("S12").Value = 0.2256036852755888
$ws.Range("T12").Value = 0.2256036852755887

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Clcf1"
$ws.Range("C13").Value = "Crlf1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 6.825289333333334
$ws.Range("H13").Value = 20.475868
$ws.Range("I13").Value = 0.5251361975000832
$ws.Range("J13").Value = 0.5251361975000832
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.06894633333333333
$ws.Range("N13").Value = 0.206839
$ws.Range("O13").Value = 0.005991935818884376
$ws.Range("P13").Value = 0.005991935818884375
$ws.Range("Q13").Value = 0.4705786734724445
$ws.Range("R13").Value = 4.235208061252
$ws.Range("S13").Value = 0.003146582391593489
$ws.Range("T13").Value = 0.003146582391593488

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Clcf1"
$ws.Range("C14").Value = "Crlf1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.099475
$ws.Range("H14").Value = 3.298425
$ws.Range("I14").Value = 0.08459335458888541
$ws.Range("J14").Value = 0.08459335458888539
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.3711636666666667
$ws.Range("N14").Value = 1.113491
$ws.Range("O14").Value = 0.0322568113697387
$ws.Range("P14").Value = 0.03225681136973869
$ws.Range("Q14").Value = 0.4080851724083333
$ws.Range("R14").Value = 3.672766551675
$ws.Range("S14").Value = 0.002728711882107097
$ws.Range("T14").Value = 0.002728711882107095

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Clcf1"
$ws.Range("C15").Value = "Crlf1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.099475
$ws.Range("H15").Value = 3.298425
$ws.Range("I15").Value = 0.08459335458888541
$ws.Range("J15").Value = 0.08459335458888539
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 6.123096
$ws.Range("N15").Value = 18.369288
$ws.Range("O15").Value = 0.5321413985496108
$ws.Range("P15").Value = 0.5321413985496107
$ws.Range("Q15").Value = 6.7321909746
$ws.Range("R15").Value = 60.5897187714
$ws.Range("S15").Value = 0.04501562601893262
$ws.Range("T15").Value = 0.0450156260189326

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Clcf1"
$ws.Range("C16").Value = "Crlf1"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.099475
$ws.Range("H16").Value = 3.298425
$ws.Range("I16").Value = 0.08459335458888541
$ws.Range("J16").Value = 0.08459335458888539
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 4.943314666666667
$ws.Range("N16").Value = 14.829944
$ws.Range("O16").Value = 0.4296098542617661
$ws.Range("P16").Value = 0.4296098542617661
$ws.Range("Q16").Value = 5.435050893133333
$ws.Range("R16").Value = 48.9154580382
$ws.Range("S16").Value = 0.03634213873644496
$ws.Range("T16").Value = 0.03634213873644496

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Clcf1"
$ws.Range("C17").Value = "Crlf1"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.099475
$ws.Range("H17").Value = 3.298425
$ws.Range("I17").Value = 0.08459335458888541
$ws.Range("J17").Value = 0.08459335458888539
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.06894633333333333
$ws.Range("N17").Value = 0.206839
$ws.Range("O17").Value = 0.005991935818884376
$ws.Range("P17").Value = 0.005991935818884375
$ws.Range("Q17").Value = 0.07580476984166666
$ws.Range("R17").Value = 0.682242928575
$ws.Range("S17").Value = 0.0005068779514007296
$ws.Range("T17").Value = 0.0005068779514007292

